$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data block lives in columns B:K, rows 2:16 (one quarterly "error"
# series per row). This edit adds a new observation to the first data
# row (row 2) and shifts every row's values one column to the left,
# dropping the oldest (leftmost) observation - i.e. a rolling window
# advanced by one period. Row 2 is the only row that still spans the
# full B:K width (it gains the new value in column K); every other row
# ends up one column shorter than before.

$newValue = -0.3663687737149753

for ($r = 2; $r -le 16; $r++) {
    # Read current row values (columns B..K => 2..11)
    $vals = @()
    for ($c = 2; $c -le 11; $c++) {
        $vals += $ws.Cells.Item($r, $c).Value2
    }

    # Drop the first (oldest) value - this is the left shift.
    $shifted = @()
    for ($i = 1; $i -lt $vals.Length; $i++) {
        $shifted += $vals[$i]
    }

    # Row 2 gets a brand-new observation appended at the end (column K).
    if ($r -eq 2) {
        $shifted += $newValue
    }

    # Write back the shifted values starting at column B, then clear
    # whatever trails (the row is now one cell shorter, except row 2).
    $c = 2
    foreach ($v in $shifted) {
        $ws.Cells.Item($r, $c).Value = $v
        $c++
    }
    for (; $c -le 11; $c++) {
        $ws.Cells.Item($r, $c).ClearContents()
    }
}
